$wb = $excel.ActiveWorkbook

# --- Update the Milestones table on "Project Tracker" with Start/End dates ---
$tracker = $wb.Worksheets.Item("Project Tracker")

# Row 7 (Position 2 - System Testing)
$tracker.Range("C7").Value = 43870
$tracker.Range("D7").Value = 43892

# Row 8 (Position 3 - Subsystem Testing)
$tracker.Range("C8").Value = 43886
$tracker.Range("D8").Value = 43892

# Row 9 (Position 4 - Unit Testing)
$tracker.Range("C9").Value = 43877
$tracker.Range("D9").Value = 43892

# Recalculate so dependent formulas/charts pick up the new values
$excel.Calculate()

# Update the selection on the "Project Tracker" sheet
$tracker.Range("C8").Select()

# --- Make "Project Chart" the active / displayed sheet ---
$chartSheet = $wb.Worksheets.Item("Project Chart")
$chartSheet.Activate()
